$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 108014880
$ws.Range("J70").Value = 55556732
$ws.Range("L70").Value = 166670196
$ws.Range("N70").Value = -166670736
$ws.Range("H73").Value = 108014880
$ws.Range("J73").Value = 55556732
$ws.Range("L73").Value = 166670196
$ws.Range("N73").Value = -166672068
$ws.Range("H92").Value = 2406.3
$ws.Range("I92").Value = 2093
$ws.Range("K92").Value = 2093
$ws.Range("M92").Value = -845
$ws.Range("H121").Value = 395.77777
$ws.Range("J121").Value = 395.25
$ws.Range("L121").Value = 1185.75
$ws.Range("N121").Value = -4679.75
$ws.Range("H132").Value = 6608.5713
$ws.Range("I132").Value = 3441.2173
$ws.Range("K132").Value = 10323.6519
$ws.Range("M132").Value = -7793.651899999999
$ws.Range("H135").Value = 3145.1
$ws.Range("J135").Value = 6856.909
$ws.Range("L135").Value = 61712.181
$ws.Range("N135").Value = -66782.181
$ws.Range("H137").Value = 20836156
$ws.Range("I137").Value = 26318422
$ws.Range("K137").Value = 78955266
$ws.Range("M137").Value = -78952716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3188.5232
$ws.Range("I32").Value = 2924.224
$ws.Range("K32").Value = 2924.224
$ws.Range("M32").Value = -2637.224
$ws.Range("H61").Value = 118335920
$ws.Range("I61").Value = 140003100
$ws.Range("J61").Value = 9999999
$ws.Range("K61").Value = 140003100
$ws.Range("L61").Value = 9999999
$ws.Range("M61").Value = -140002888
$ws.Range("N61").Value = -10000423
$ws.Range("H74").Value = 2738.7334
$ws.Range("I74").Value = 2927
$ws.Range("K74").Value = 2927
$ws.Range("M74").Value = -2053
$ws.Range("H77").Value = 2738.7334
$ws.Range("I77").Value = 2927
$ws.Range("K77").Value = 14635
$ws.Range("M77").Value = -10267
$ws.Range("H88").Value = 4163.6924
$ws.Range("I88").Value = 3696.5
$ws.Range("J88").Value = 4371.3335
$ws.Range("K88").Value = 3696.5
$ws.Range("L88").Value = 4371.3335
$ws.Range("M88").Value = -3290.5
$ws.Range("N88").Value = -5183.3335
$ws.Range("H91").Value = 4163.6924
$ws.Range("I91").Value = 3696.5
$ws.Range("J91").Value = 4371.3335
$ws.Range("K91").Value = 3696.5
$ws.Range("L91").Value = 4371.3335
$ws.Range("M91").Value = -2292.5
$ws.Range("N91").Value = -7179.3335
$ws.Range("H102").Value = 2657.9092
$ws.Range("I102").Value = 1804.3334
$ws.Range("K102").Value = 1804.3334
$ws.Range("M102").Value = -182.3334
$ws.Range("H132").Value = 10641889
$ws.Range("I132").Value = 3467.7896
$ws.Range("J132").Value = 55559664
$ws.Range("K132").Value = 10403.3688
$ws.Range("L132").Value = 166678992
$ws.Range("M132").Value = -7873.3688
$ws.Range("N132").Value = -166684052
$ws.Range("H136").Value = 118335920
$ws.Range("I136").Value = 140003100
$ws.Range("J136").Value = 9999999
$ws.Range("K136").Value = 420009300
$ws.Range("L136").Value = 29999997
$ws.Range("M136").Value = -420006750
$ws.Range("N136").Value = -30005097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 732676.0600000001
$ws.Range("I86").Value = 1298169.9
$ws.Range("K86").Value = 1298169.9
$ws.Range("M86").Value = -1297046.9
$ws.Range("H89").Value = 732676.0600000001
$ws.Range("I89").Value = 1298169.9
$ws.Range("K89").Value = 6490849.5
$ws.Range("M89").Value = -6485233.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15627662
$ws.Range("I31").Value = 31253124
$ws.Range("J31").Value = 2200.9375
$ws.Range("K31").Value = 31253124
$ws.Range("L31").Value = 2200.9375
$ws.Range("M31").Value = -31252829
$ws.Range("N31").Value = -2790.9375
$ws.Range("H34").Value = 15627662
$ws.Range("I34").Value = 31253124
$ws.Range("J34").Value = 2200.9375
$ws.Range("K34").Value = 31253124
$ws.Range("L34").Value = 2200.9375
$ws.Range("M34").Value = -31252922
$ws.Range("N34").Value = -2604.9375
$ws.Range("H58").Value = 2894.6924
$ws.Range("I58").Value = 1843.25
$ws.Range("K58").Value = 1843.25
$ws.Range("M58").Value = -1640.25
$ws.Range("H99").Value = 26979.111
$ws.Range("I99").Value = 8497.333000000001
$ws.Range("K99").Value = 8497.333000000001
$ws.Range("M99").Value = -6999.333000000001
$ws.Range("H126").Value = 26979.111
$ws.Range("I126").Value = 8497.333000000001
$ws.Range("K126").Value = 25491.999
$ws.Range("M126").Value = -23021.999
$ws.Range("H132").Value = 2014.8462
$ws.Range("I132").Value = 2014.8462
$ws.Range("K132").Value = 6044.5386
$ws.Range("M132").Value = -3514.5386
$ws.Range("H136").Value = 2894.6924
$ws.Range("I136").Value = 1843.25
$ws.Range("K136").Value = 5529.75
$ws.Range("M136").Value = -2979.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 1696
$ws.Range("J5").Value = 1505.5
$ws.Range("K5").Value = 5088
$ws.Range("L5").Value = 4516.5
$ws.Range("M5").Value = -4976
$ws.Range("N5").Value = -4740.5
$ws.Range("H63").Value = 30665.2
$ws.Range("J63").Value = 30833
$ws.Range("L63").Value = 92499
$ws.Range("N63").Value = -93997
$ws.Range("H66").Value = 30665.2
$ws.Range("J66").Value = 30833
$ws.Range("L66").Value = 277497
$ws.Range("N66").Value = -284985
$ws.Range("H107").Value = 3960560.2
$ws.Range("I107").Value = 2569.1667
$ws.Range("K107").Value = 7707.500100000001
$ws.Range("M107").Value = -5787.500100000001
$ws.Range("H121").Value = 7031.923
$ws.Range("J121").Value = 9231.223
$ws.Range("L121").Value = 27693.669
$ws.Range("N121").Value = -30313.669
$ws.Range("H132").Value = 1883
$ws.Range("I132").Value = 1825.75
$ws.Range("J132").Value = 1997.5
$ws.Range("K132").Value = 16431.75
$ws.Range("L132").Value = 17977.5
$ws.Range("M132").Value = -13901.75
$ws.Range("N132").Value = -23037.5
$ws.Range("I135").Value = 1696
$ws.Range("J135").Value = 1505.5
$ws.Range("K135").Value = 15264
$ws.Range("L135").Value = 13549.5
$ws.Range("M135").Value = -12729
$ws.Range("N135").Value = -18619.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 99253.336
$ws.Range("J119").Value = 99253.336
$ws.Range("L119").Value = 99253.336
$ws.Range("N119").Value = -108929.336
$ws.Range("H122").Value = 3829.1667
$ws.Range("J122").Value = 4999
$ws.Range("L122").Value = 14997
$ws.Range("N122").Value = -19897
$ws.Range("H132").Value = 2096.923
$ws.Range("I132").Value = 2121.75
$ws.Range("J132").Value = 1799
$ws.Range("K132").Value = 6365.25
$ws.Range("L132").Value = 5397
$ws.Range("M132").Value = -3835.25
$ws.Range("N132").Value = -10457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 27000
$ws.Range("I23").Value = 27000
$ws.Range("K23").Value = 27000
$ws.Range("M23").Value = -26770
$ws.Range("H55").Value = 927.3077
$ws.Range("J55").Value = 1145.909
$ws.Range("L55").Value = 1145.909
$ws.Range("N55").Value = -1491.909
$ws.Range("H68").Value = 2606652.5
$ws.Range("J68").Value = 2473.75
$ws.Range("L68").Value = 2473.75
$ws.Range("N68").Value = -3971.75
$ws.Range("H71").Value = 2606652.5
$ws.Range("J71").Value = 2473.75
$ws.Range("L71").Value = 12368.75
$ws.Range("N71").Value = -19856.75
$ws.Range("H100").Value = 20858202
$ws.Range("I100").Value = 3695.6
$ws.Range("J100").Value = 35754280
$ws.Range("K100").Value = 3695.6
$ws.Range("L100").Value = 35754280
$ws.Range("M100").Value = -3154.6
$ws.Range("N100").Value = -35755362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3421.8333
$ws.Range("I107").Value = 1319.1052
$ws.Range("K107").Value = 3957.3156
$ws.Range("M107").Value = -2037.3156
$ws.Range("H113").Value = 753.6667
$ws.Range("I113").Value = 772.5
$ws.Range("J113").Value = 603
$ws.Range("K113").Value = 2317.5
$ws.Range("L113").Value = 1809
$ws.Range("M113").Value = -147.5
$ws.Range("N113").Value = -6149
$ws.Range("H122").Value = 2795.2646
$ws.Range("I122").Value = 2751.25
$ws.Range("K122").Value = 8253.75
$ws.Range("M122").Value = -5803.75
$ws.Range("H135").Value = 84825
$ws.Range("J135").Value = 119000
$ws.Range("L135").Value = 119000
$ws.Range("N135").Value = -129140
